# Apply the edit described by the diff:
#  - A1 date value changes from 45406 (24-Apr-2024) to 45436 (24-May-2024)
#  - D19 price value changes from 440 to 1010
#  - the <mergeCells> list is rewritten in a new order (same 8 ranges)
#
# Re-ordering the merge list requires rebuilding the merges (unmerge, then
# re-merge in the desired sequence). Because re-merging restamps formatting
# on the covered (non-anchor) cells of each merged range, we snapshot their
# original formatting first and restore it afterwards so the only visible
# changes are the two cell values and the merge order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$xlPasteFormats = -4122

# --- value edits ---------------------------------------------------------
$ws.Range("A1").Value = 45436
$ws.Range("D19").Value = 1010

# --- snapshot formatting of every cell inside a merged range -------------
# (re-merging restamps formatting/borders on ALL cells in the range,
# including the anchor cell, so every one of them needs a backup)
# (use an unused scratch area far below the sheet's real data as scratch pad)
$ws.Range("A12:D12").Copy()
$ws.Range("A212:D212").PasteSpecial($xlPasteFormats)
$ws.Range("A13:D13").Copy()
$ws.Range("A213:D213").PasteSpecial($xlPasteFormats)
$ws.Range("A14:D14").Copy()
$ws.Range("A214:D214").PasteSpecial($xlPasteFormats)
$ws.Range("B17:C17").Copy()
$ws.Range("B217:C217").PasteSpecial($xlPasteFormats)
$ws.Range("B18:C18").Copy()
$ws.Range("B218:C218").PasteSpecial($xlPasteFormats)
$ws.Range("B19:C19").Copy()
$ws.Range("B219:C219").PasteSpecial($xlPasteFormats)
$ws.Range("A11").Copy()
$ws.Range("A211").PasteSpecial($xlPasteFormats)
$ws.Range("A1").Copy()
$ws.Range("A201").PasteSpecial($xlPasteFormats)

# --- unmerge all 8 merged ranges ------------------------------------------
$order = @("A1:D1", "A12:D12", "B19:C19", "A13:D13", "A11:D11", "B18:C18", "B17:C17", "A14:D14")
foreach ($ref in $order) {
    $ws.Range($ref).UnMerge()
}

# --- re-merge in the new target order -------------------------------------
foreach ($ref in $order) {
    $ws.Range($ref).Merge()
}

# --- restore formatting of every merged cell from the snapshot -----------
$ws.Range("A212:D212").Copy()
$ws.Range("A12:D12").PasteSpecial($xlPasteFormats)
$ws.Range("A213:D213").Copy()
$ws.Range("A13:D13").PasteSpecial($xlPasteFormats)
$ws.Range("A214:D214").Copy()
$ws.Range("A14:D14").PasteSpecial($xlPasteFormats)
$ws.Range("B217:C217").Copy()
$ws.Range("B17:C17").PasteSpecial($xlPasteFormats)
$ws.Range("B218:C218").Copy()
$ws.Range("B18:C18").PasteSpecial($xlPasteFormats)
$ws.Range("B219:C219").Copy()
$ws.Range("B19:C19").PasteSpecial($xlPasteFormats)
$ws.Range("A211").Copy()
$ws.Range("A11").PasteSpecial($xlPasteFormats)
$ws.Range("A201").Copy()
$ws.Range("A1").PasteSpecial($xlPasteFormats)

# A1:D1 and A11:D11 never had their covered cells (B/C/D) materialized in
# the workbook; reset them back to the default style so they stay absent.
$ws.Range("B1:D1").Style = "Normal"
$ws.Range("B11:D11").Style = "Normal"

# --- clean up scratch area --------------------------------------------------
$ws.Range("A201:D219").Clear()
